$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name
$ws.Range("C3").Value = "Meng Wang"

# Row 7 - __init__ / Attribute set to input value (happy path)
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = @"
name = "Intermediate Software Development" 
department = COMPUTER_SCIENCE
credit_hours = 6
"@
$ws.Range("G7").Value = "Client object created with expected attributes value based on method inputs"

# Row 8 - __init__ / Exception raised when name is blank
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = @"
name = " " 
department = COMPUTER_SCIENCE
credit_hours = 6
"@
$ws.Range("G8").Value = "ValueError"

# Row 9 - __init__ / Exception raised when invalid department
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = @"
name = "Intermediate Software Development" 
department = INVALID
credit_hours = 6
"@
$ws.Range("G9").Value = "ValueError"

# Row 10 - __init__ / Exception raised when non-numeric credit_hours
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = @"
name = "Intermediate Software Development" 
department = COMPUTER_SCIENCE
credit_hours = "six"
"@
$ws.Range("G10").Value = "ValueError"

# Row 11 - name / returns name attribute
$ws.Range("E11").Value = @"
Course("Intermediate Software Development",
Department.COMPUTER_SCIENCE,6)
"@
$ws.Range("F11").Value = "None"
$ws.Range("G11").Value = 'course._Course__name = "Intermeddiate Software Development"'

# Row 12 - department / returns department attribute
$ws.Range("E12").Value = @"
Course("Intermediate Software Development",
Department.COMPUTER_SCIENCE,6)
"@
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = "Course._Course__department = Department.COMPUTER_SCIENCE"

# Row 13 - credit_hours / returns credit_hours attribute
$ws.Range("E13").Value = @"
Course("Intermediate Software Development",
Department.COMPUTER_SCIENCE,6)
"@
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = "course._Cource__credit_hours = 6"

# Row 14 - __str__ / returns string in expected format
$ws.Range("E14").Value = @"
Course("Intermediate Software Development",
Department.COMPUTER_SCIENCE,6)
"@
$ws.Range("F14").Value = "None"
$ws.Range("G14").Value = @"
"Course: Intermediate
Software Development\n
Department: Computer
Science\n
Credit Hours:6"
"@

# Move selection / scroll position to roughly match the saved view
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G14").Select()
